# Applies the "Add performance test references" commit.
#
# Summary of changes (see XML diff):
#  - Slide 10 ("Example: Manually Create an Overlay Network", text-only):
#      * merge the "ocker Daemon, and set IP ranges " / "with " / "--" runs
#      * merge the " devices and " / "routes with iproute2 and bridge" runs
#      * merge the "Create ARP entries across all " / "hosts with iproute2" runs
#  - Slide 13 ("VXLAN backend"):
#      * add a new paragraph "Performance test refers to a post of mine"
#        with "a post of mine" hyperlinked
#  - Slide 3 ("outline"):
#      * merge " backend in " / "flannel" runs
#      * merge "Brief introduction to " / "implementation of " runs
#  - Slide 9 ("Example: Manually Create an Overlay Network", with picture):
#      * merge the " Create " / "an " / "Overlay Network" title runs

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10: collapse several same-formatted run splits back into single runs.
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$body10 = $s10.Shapes.Item(2).TextFrame.TextRange

$para = $body10.Paragraphs(1)
$sub = $para.Characters(20, 39)
$sub.Text = "ocker Daemon, and set IP ranges with --"

$para = $body10.Paragraphs(2)
$sub = $para.Characters(38, 44)
$sub.Text = " devices and routes with iproute2 and bridge"

$para = $body10.Paragraphs(3)
$sub = $para.Characters(1, 49)
$sub.Text = "Create ARP entries across all hosts with iproute2"

# ---------------------------------------------------------------------------
# Slide 13: split the last sentence and add the new "performance test" line,
# linking "a post of mine" to the author's blog post.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$body13 = $s13.Shapes.Item(2).TextFrame.TextRange

$body13.InsertAfter([char]13 + "Performance test refers to a post of mine")

$full13 = $body13.Text
$idx = $full13.IndexOf("a post of mine")
$linkRange = $body13.Characters($idx + 1, "a post of mine".Length)
$linkRange.ActionSettings.Item(1).Action = 7
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://ieevee.com/tech/2017/03/24/docker-vxlan-network.html"

# ---------------------------------------------------------------------------
# Slide 3: collapse outline-bullet run splits back into single runs.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

$para = $body3.Paragraphs(4)
$sub = $para.Characters(6, 19)
$sub.Text = " backend in flannel"

$para = $body3.Paragraphs(5)
$sub = $para.Characters(1, 40)
$sub.Text = "Brief introduction to implementation of "

# ---------------------------------------------------------------------------
# Slide 9: collapse title run splits back into a single run.
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Item(1).TextFrame.TextRange

$sub = $title9.Characters(18, 26)
$sub.Text = " Create an Overlay Network"
